$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = 0.75
$ws.Range("D15").Value = 1
$ws.Range("D19").Value = 0.75
$ws.Range("D21").Value = 6.25
